# Replace the arithmetic problems in the practice-sheet table.
# Each "old" value is unique within the document, so a straightforward
# Find/Execute with Replace:=2 (wdReplaceAll) for each pair is safe and
# unambiguous.

$d = $word.ActiveDocument

$pairs = @(
    @{ Old = "392×9="; New = "471×9=" },
    @{ Old = "174×5="; New = "182×4=" },
    @{ Old = "424×3="; New = "676×8=" },
    @{ Old = "802×6="; New = "396×5=" },
    @{ Old = "514×6="; New = "401×5=" },
    @{ Old = "904×2="; New = "169×2=" },
    @{ Old = "518×9="; New = "303×2=" },
    @{ Old = "139×5="; New = "971×4=" },
    @{ Old = "431×2="; New = "113×4=" },
    @{ Old = "299×2="; New = "665×8=" },
    @{ Old = "333×4="; New = "885×8=" },
    @{ Old = "422×2="; New = "600×2=" },
    @{ Old = "115×9="; New = "441×3=" },
    @{ Old = "338×6="; New = "877×7=" },
    @{ Old = "466×3="; New = "179×4=" },
    @{ Old = "985×5="; New = "145×4=" },
    @{ Old = "975×5="; New = "845×8=" },
    @{ Old = "298×3="; New = "944×5=" },
    @{ Old = "238×9="; New = "743×8=" },
    @{ Old = "506×5="; New = "180×7=" },
    @{ Old = "701×8="; New = "565×5=" },
    @{ Old = "204×9="; New = "215×9=" },
    @{ Old = "437×9="; New = "461×6=" },
    @{ Old = "861×4="; New = "558×6=" },
    @{ Old = "985×2="; New = "343×6=" }
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
